# Add a new worksheet "Thicknesses" containing sample thickness data,
# make it the active sheet, and update the selection on Sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Create the new sheet after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Thicknesses"

# Headers
$ws2.Range("A1").Value = "Sample"
$ws2.Range("B1").Value = "Thickness [A]"

# Data rows: Sample # -> Thickness [A]
$data = @(
    @(2, 233),
    @(15, 205.75),
    @(21, 237.5),
    @(22, 282.25),
    @(24, 304.25),
    @(39, 257),
    @(50, 290.25),
    @(53, 234.5),
    @(58, 321.5),
    @(81, 265.25),
    @(1, 471.33333),
    @(3, 441.33333),
    @(12, 447.33333),
    @(23, 167),
    @(61, 375.66667),
    @(80, 300)
)

$row = 2
foreach ($pair in $data) {
    $ws2.Cells.Item($row, 1).Value = $pair[0]
    $ws2.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Bold the sample numbers in column A (rows 2-17)
$ws2.Range("A2:A17").Font.Bold = $true

# Selections / active sheet
$ws1.Range("C3:C18").Select()
$ws2.Activate()
$ws2.Range("E14").Select()
